$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.83"
$ws.Range("E2").Value = "'-4.09%"
$ws.Range("D3").Value = "'54.40"
$ws.Range("E3").Value = "'10.54%"
$ws.Range("D4").Value = "'5.074"
$ws.Range("E4").Value = "'-5.22%"
$ws.Range("D5").Value = "'0.07861"
$ws.Range("E5").Value = "'-2.47%"
$ws.Range("D6").Value = "'4.529"
$ws.Range("E6").Value = "'-1.67%"
$ws.Range("D7").Value = "'1.375"
$ws.Range("E7").Value = "'-1.83%"
$ws.Range("D8").Value = "'1.728"
$ws.Range("E8").Value = "'5.39%"
$ws.Range("D9").Value = "'0.1240"
$ws.Range("E9").Value = "'-4.02%"
$ws.Range("D10").Value = "'0.2003"
$ws.Range("E10").Value = "'1.89%"
$ws.Range("D11").Value = "'0.04710"
$ws.Range("E11").Value = "'0.88%"
$ws.Range("D12").Value = "'0.09452"
$ws.Range("E12").Value = "'-0.70%"
$ws.Range("D13").Value = "'0.1041"
$ws.Range("E13").Value = "'-0.23%"
$ws.Range("D14").Value = "'0.001256"
$ws.Range("E14").Value = "'-4.57%"
$ws.Range("D15").Value = "'0.005770"
$ws.Range("E15").Value = "'-0.22%"
$ws.Range("E16").Value = "'2,021.91%"
$ws.Range("E17").Value = "'-0.57%"
$ws.Range("D18").Value = "'2.413"
$ws.Range("E18").Value = "'-1.10%"
$ws.Range("E19").Value = "'-2.61%"
$ws.Range("D20").Value = "'8.002"
$ws.Range("E20").Value = "'-0.49%"
$ws.Range("D21").Value = "'0.1360"
$ws.Range("E21").Value = "'-1.05%"
$ws.Range("D23").Value = "'0.04159"
$ws.Range("E23").Value = "'-0.73%"
$ws.Range("D24").Value = "'0.001261"
$ws.Range("E24").Value = "'-4.05%"
$ws.Range("D25").Value = "'0.003929"
$ws.Range("E25").Value = "'-8.73%"
$ws.Range("D26").Value = "'0.0001350"
$ws.Range("E26").Value = "'0.32%"
$ws.Range("D38").Value = "'0.02610"
$ws.Range("E38").Value = "'-4.95%"
$ws.Range("D39").Value = "'0.05871"
$ws.Range("E39").Value = "'-7.43%"
$ws.Range("D40").Value = "'0.009849"
$ws.Range("E40").Value = "'-5.72%"
$ws.Range("D41").Value = "'0.007959"
$ws.Range("E41").Value = "'-0.97%"
$ws.Range("D42").Value = "'0.1440"
$ws.Range("E42").Value = "'-1.46%"
$ws.Range("D43").Value = "'0.008201"
$ws.Range("E43").Value = "'3.99%"
$ws.Range("D44").Value = "'0.008361"
$ws.Range("E44").Value = "'-3.05%"
$ws.Range("D45").Value = "'0.3360"
$ws.Range("E45").Value = "'-4.17%"
$ws.Range("D46").Value = "'0.00007327"
$ws.Range("E46").Value = "'10.49%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.27%"
$ws.Range("D48").Value = "'0.05680"
$ws.Range("E48").Value = "'1.60%"
$ws.Range("D49").Value = "'0.002613"
$ws.Range("E49").Value = "'-34.48%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'0.27%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'0.27%"
